# Applies:
#  1) Table style change on the table (graphicFrame) in slide 16:
#     {49011210-8EBF-4C54-AC75-F3183683E9F6} (custom "Table_0") ->
#     {82EF40D3-B3A4-419D-96AD-E0976F85F8ED} (built-in table style)
#  2) The presentation's active theme color scheme (theme used by the
#     slide master / all slides) is switched from the "Integral" theme
#     palette back to the default "Office" theme palette.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$targetSlideIndex = 16
$newTableStyleId   = "{82EF40D3-B3A4-419D-96AD-E0976F85F8ED}"

$slide = $p.Slides.Item($targetSlideIndex)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle($newTableStyleId)
    }
}

# --- 2) Theme colors: Integral -> Office -----------------------------
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. Values are packed BGR integers (as COM RGB() would
# produce: R + G*256 + B*65536) for the stock "Office" theme palette.
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeRGB[$i - 1]
}
